$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates parsed from the target diff: column D (Price), E (Volume 1h)
# and G (Hora) for each data row. Rows whose Price/Volume are placeholders ("--")
# in the source only get their G (Hora) value bumped from 2 to 3.
$updates = @(
    @{ Row = 2; D = '286.89'; E = '4.03%'; G = '3' },
    @{ Row = 3; D = '28.37'; E = '4.40%'; G = '3' },
    @{ Row = 4; D = '4.933'; E = '1.60%'; G = '3' },
    @{ Row = 5; D = '0.06555'; E = '2.40%'; G = '3' },
    @{ Row = 6; D = '7.260'; E = '4.82%'; G = '3' },
    @{ Row = 7; D = '1.347'; E = '10.23%'; G = '3' },
    @{ Row = 8; D = '0.9170'; E = '4.14%'; G = '3' },
    @{ Row = 9; D = '0.1566'; E = '3.47%'; G = '3' },
    @{ Row = 10; D = '0.06409'; E = '25.91%'; G = '3' },
    @{ Row = 11; D = '0.07684'; E = '0.91%'; G = '3' },
    @{ Row = 12; D = '0.02978'; E = '0.02%'; G = '3' },
    @{ Row = 13; D = '0.08972'; E = '-0.35%'; G = '3' },
    @{ Row = 14; D = '0.001607'; E = '2.54%'; G = '3' },
    @{ Row = 15; D = '0.0006567'; E = '1.98%'; G = '3' },
    @{ Row = 16; D = '0.005998'; E = '-3.19%'; G = '3' },
    @{ Row = 17; D = '3.484'; E = '0.53%'; G = '3' },
    @{ Row = 18; D = '3.389'; E = '2.47%'; G = '3' },
    @{ Row = 19; D = '2.237'; E = '-2.07%'; G = '3' },
    @{ Row = 20; G = '3' },
    @{ Row = 21; E = '-0.47%'; G = '3' },
    @{ Row = 22; D = '3.996'; E = '1.94%'; G = '3' },
    @{ Row = 23; D = '0.04468'; E = '0.93%'; G = '3' },
    @{ Row = 24; D = '0.1521'; E = '10.24%'; G = '3' },
    @{ Row = 25; E = '0.70%'; G = '3' },
    @{ Row = 26; D = '0.004348'; E = '1.84%'; G = '3' },
    @{ Row = 27; G = '3' },
    @{ Row = 28; D = '0.0001179'; E = '-1.92%'; G = '3' },
    @{ Row = 29; E = '-15.72%'; G = '3' },
    @{ Row = 30; G = '3' },
    @{ Row = 31; G = '3' },
    @{ Row = 32; G = '3' },
    @{ Row = 33; G = '3' },
    @{ Row = 34; G = '3' },
    @{ Row = 35; G = '3' },
    @{ Row = 36; G = '3' },
    @{ Row = 37; G = '3' },
    @{ Row = 38; G = '3' },
    @{ Row = 39; G = '3' },
    @{ Row = 40; D = '0.04148'; E = '0.24%'; G = '3' },
    @{ Row = 41; D = '0.006881'; E = '0.73%'; G = '3' },
    @{ Row = 42; D = '0.1412'; E = '20.30%'; G = '3' },
    @{ Row = 43; D = '0.002049'; E = '-4.90%'; G = '3' },
    @{ Row = 44; D = '0.01252'; E = '5.63%'; G = '3' },
    @{ Row = 45; D = '0.00005551'; E = '7.14%'; G = '3' },
    @{ Row = 46; D = '1.561'; E = '-5.55%'; G = '3' },
    @{ Row = 47; D = '0.01849'; E = '-7.76%'; G = '3' },
    @{ Row = 48; G = '3' },
    @{ Row = 49; G = '3' },
    @{ Row = 50; G = '3' },
    @{ Row = 51; G = '3' }
)

foreach ($u in $updates) {
    foreach ($col in @("D", "E", "G")) {
        if ($u.ContainsKey($col)) {
            $ref = "$col$($u.Row)"
            # Force text storage so numeric-looking strings (e.g. "286.89", "3")
            # are not auto-coerced into Number cells, matching the source's
            # inline-string (text) cell type.
            $ws.Range($ref).NumberFormat = "@"
            $ws.Range($ref).Value = $u[$col]
            $ws.Range($ref).Style = "Normal"
        }
    }
}